$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Structural changes: insert new column L (shifts old L..P to M..Q),
#        and append 5 new rows at the bottom (17..21) to make room for the
#        5 new data points that were added to the experiment table.
$ws.Columns("L").Insert()
$ws.Rows("17:21").Insert()

# --- 2. Write every target cell's final content (values / formulas).
#        New references are written before any stale cell is cleared so
#        that shared-string table ordering/indices never transiently drop
#        a string that is still needed (avoids spurious re-numbering).
$ws.Range("A1").Value = "payload"
$ws.Range("B1").Value = "from baseline"
$ws.Range("C1").Value = "hidden"
$ws.Range("D1").Value = "delta-from-defaults"
$ws.Range("E1").Value = "ordering by dist"
$ws.Range("F1").Value = "order by diffcount"
$ws.Range("G1").Value = "order by worst"
$ws.Range("H1").Value = "orderby compr"
$ws.Range("I1").Value = "order by`nworst compr"
$ws.Range("J1").Value = "interleave16"
$ws.Range("K1").Value = "interleave8"
$ws.Range("L1").Value = "delta16sign"
$ws.Range("M1").Value = "delta-16high"
$ws.Range("N1").Value = "16low"
$ws.Range("O1").Value = "8high"
$ws.Range("P1").Value = "8low"
$ws.Range("Q1").Value = "Note"
$ws.Range("A2").Value = 12883
$ws.Range("B2").Formula = "=A2-13052"
$ws.Range("D2").Value = "x"
$ws.Range("E2").Value = "x"
$ws.Range("Q2").Value = "Reordering does have impact"
$ws.Range("A3").Value = 12887
$ws.Range("B3").Formula = "=A3-13052"
$ws.Range("D3").Value = "x"
$ws.Range("E3").Value = "x"
$ws.Range("J3").Value = "x"
$ws.Range("Q3").Value = "interleaving params is interesting but not helping."
$ws.Range("A4").Value = 12908
$ws.Range("B4").Formula = "=A4-13052"
$ws.Range("D4").Value = "x"
$ws.Range("F4").Value = "x"
$ws.Range("Q4").Value = "instead of total distance, count # of bytes that would change. Slightly worse."
$ws.Range("A5").Value = 12909
$ws.Range("B5").Formula = "=A5-13052"
$ws.Range("D5").Value = "x"
$ws.Range("E5").Value = "x"
$ws.Range("L5").Value = "x"
$ws.Range("Q5").Value = "ok the awful results earlier wsa from not handling the sign correctly."
$ws.Range("A6").Value = 12911
$ws.Range("B6").Formula = "=A6-13052"
$ws.Range("D6").Value = "x"
$ws.Range("E6").Value = "x"
$ws.Range("J6").Value = "x"
$ws.Range("L6").Value = "x"
$ws.Range("A7").Value = 12916
$ws.Range("B7").Formula = "=A7-13052"
$ws.Range("D7").Value = "x"
$ws.Range("E7").Value = "x"
$ws.Range("K7").Value = "x"
$ws.Range("Q7").Value = "interleaving by byte is also not helpful"
$ws.Range("A8").Value = 12920
$ws.Range("B8").Formula = "=A8-13052"
$ws.Range("D8").Value = "x"
$ws.Range("H8").Value = "x"
$ws.Range("L8").Value = "x"
$ws.Range("A9").Value = 12934
$ws.Range("B9").Formula = "=A9-13052"
$ws.Range("D9").Value = "x"
$ws.Range("H9").Value = "x"
$ws.Range("Q9").Value = "how does this possibly perform worse than dist? Well it does."
$ws.Range("A10").Value = 13043
$ws.Range("B10").Formula = "=A10-13052"
$ws.Range("D10").Value = "x"
$ws.Range("G10").Value = "x"
$ws.Range("A11").Value = 13052
$ws.Range("B11").Formula = "=A11-13052"
$ws.Range("D11").Value = "x"
$ws.Range("A12").Value = 13078
$ws.Range("B12").Formula = "=A12-13052"
$ws.Range("D12").Value = "x"
$ws.Range("I12").Value = "x"
$ws.Range("A13").Value = 13086
$ws.Range("B13").Formula = "=A13-13052"
$ws.Range("E13").Value = "x"
$ws.Range("A14").Value = 13093
$ws.Range("B14").Formula = "=A14-13052"
$ws.Range("H14").Value = "x"
$ws.Range("J14").Value = "x"
$ws.Range("L14").Value = "x"
$ws.Range("A15").Value = 13120
$ws.Range("B15").Formula = "=A15-13052"
$ws.Range("H15").Value = "x"
$ws.Range("L15").Value = "x"
$ws.Range("A16").Value = 13127
$ws.Range("B16").Formula = "=A16-13052"
$ws.Range("E16").Value = "x"
$ws.Range("J16").Value = "x"
$ws.Range("A17").Value = 13231
$ws.Range("B17").Formula = "=A17-13052"
$ws.Range("Q17").Value = "Expected more impact from NO delta encoding"
$ws.Range("A18").Value = 13241
$ws.Range("B18").Formula = "=A18-13052"
$ws.Range("G18").Value = "x"
$ws.Range("Q18").Value = "even without ANY deltas, and WORST ordering, quite low impact."
$ws.Range("A19").Formula = "=14000"
$ws.Range("B19").Formula = "=A19-13052"
$ws.Range("C19").Value = 36639
$ws.Range("E19").Value = "x"
$ws.Range("M19").Value = "x"
$ws.Range("Q19").Value = "astonishing that ADDING a delta encoding would bloat so much, WITH reordering"
$ws.Range("A20").Formula = "=14000"
$ws.Range("B20").Formula = "=A20-13052"
$ws.Range("C20").Value = 32753
$ws.Range("E20").Value = "x"
$ws.Range("O20").Value = "x"
$ws.Range("Q20").Value = "delta encoding on byte basis is best"
$ws.Range("A21").Formula = "=14000"
$ws.Range("B21").Formula = "=A21-13052"
$ws.Range("C21").Value = 26499
$ws.Range("E21").Value = "x"
$ws.Range("P21").Value = "x"
$ws.Range("Q21").Value = "Putting the sign on the low bit does appear to help but it doesn't matter because this is awful"
# --- 3. Remove the now-stale cell contents left over at old positions
#        that are not part of the final layout (safe now: every shared
#        string above has already been re-written at its new home).
$ws.Range("K5").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("Q6").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("Q11").ClearContents()
$ws.Range("G12").ClearContents()
$ws.Range("Q12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("Q13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("O14").ClearContents()
$ws.Range("Q14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("P15").ClearContents()
$ws.Range("Q15").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("I16").ClearContents()
# --- 4. Misc view/selection tweaks to match the edited workbook.
$ws.Range("D2:D12").Select()
$excel.ActiveCell = $ws.Range("D12")
